$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.4214143333333333
$ws.Range("N2").Value = 1.264243
$ws.Range("O2").Value = 0.02434128610922473
$ws.Range("P2").Value = 0.02434128610922473
$ws.Range("Q2").Value = 0.03003307576511111
$ws.Range("R2").Value = 0.270297681886
$ws.Range("S2").Value = 0.02434128610922473
$ws.Range("T2").Value = 0.02434128610922473

# Row 3
$ws.Range("N3").Value = 36.386704
$ws.Range("O3").Value = 0.7005766871049885
$ws.Range("P3").Value = 0.7005766871049887
$ws.Range("Q3").Value = 0.8643944542897778
$ws.Range("R3").Value = 7.779550088608
$ws.Range("S3").Value = 0.7005766871049885
$ws.Range("T3").Value = 0.7005766871049887

# Row 4
$ws.Range("M4").Value = 4.762423333333333
$ws.Range("N4").Value = 14.28727
$ws.Range("O4").Value = 0.2750820267857866
$ws.Range("P4").Value = 0.2750820267857866
$ws.Range("Q4").Value = 0.3394052111711111
$ws.Range("R4").Value = 3.05464690054
$ws.Range("S4").Value = 0.2750820267857866
$ws.Range("T4").Value = 0.2750820267857866
